$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tabela1")

# New daily COVID data rows (2020-05-25 .. 2020-05-29), appended to the table.
$data = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
    @(43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0)
)

foreach ($values in $data) {
    # Grow the table by one row - this keeps the table/autoFilter ref and the
    # sheet dimension in sync with the new data automatically.
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range.Row

    for ($c = 1; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Match the number format of the row above so the freshly-grown
        # table row keeps behaving as numbers/dates instead of text.
        $cell.NumberFormat = $ws.Cells.Item($r - 1, $c).NumberFormat
        $cell.Value = $values[$c - 1]
    }
}

$lastRow = $lo.Range.Rows.Count
$ws.Range("A" + $lastRow + ":J" + $lastRow).Select()
